$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "128"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "340000.00"
$ws.Range("D12").Style = "Normal"

$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "16"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "42500.00"
$ws.Range("D32").Style = "Normal"

$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "30"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "92000.00"
$ws.Range("D33").Style = "Normal"

$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "163"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "457408.00"
$ws.Range("D35").Style = "Normal"

$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "347"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1326411.70"
$ws.Range("D37").Style = "Normal"

$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "16"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "41170.00"
$ws.Range("D40").Style = "Normal"

$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "27"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "66971.00"
$ws.Range("D43").Style = "Normal"

$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "40"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "116034.54"
$ws.Range("D45").Style = "Normal"

$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "15"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "50500.00"
$ws.Range("D47").Style = "Normal"

$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "107"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "275468.33"
$ws.Range("D50").Style = "Normal"

$ws.Range("C78").NumberFormat = "@"
$ws.Range("C78").Value = "189"
$ws.Range("C78").Style = "Normal"
$ws.Range("D78").NumberFormat = "@"
$ws.Range("D78").Value = "520693.00"
$ws.Range("D78").Style = "Normal"

$ws.Range("C79").NumberFormat = "@"
$ws.Range("C79").Value = "22"
$ws.Range("C79").Style = "Normal"
$ws.Range("D79").NumberFormat = "@"
$ws.Range("D79").Value = "65991.00"
$ws.Range("D79").Style = "Normal"

$ws.Range("C80").NumberFormat = "@"
$ws.Range("C80").Value = "433"
$ws.Range("C80").Style = "Normal"
$ws.Range("D80").NumberFormat = "@"
$ws.Range("D80").Value = "1740694.99"
$ws.Range("D80").Style = "Normal"

$ws.Range("C88").NumberFormat = "@"
$ws.Range("C88").Value = "60"
$ws.Range("C88").Style = "Normal"
$ws.Range("D88").NumberFormat = "@"
$ws.Range("D88").Value = "258736.08"
$ws.Range("D88").Style = "Normal"

$ws.Range("C90").NumberFormat = "@"
$ws.Range("C90").Value = "23"
$ws.Range("C90").Style = "Normal"
$ws.Range("D90").NumberFormat = "@"
$ws.Range("D90").Value = "49555.00"
$ws.Range("D90").Style = "Normal"

$ws.Range("C91").NumberFormat = "@"
$ws.Range("C91").Value = "55"
$ws.Range("C91").Style = "Normal"
$ws.Range("D91").NumberFormat = "@"
$ws.Range("D91").Value = "158957.14"
$ws.Range("D91").Style = "Normal"

$ws.Range("C92").NumberFormat = "@"
$ws.Range("C92").Value = "55"
$ws.Range("C92").Style = "Normal"
$ws.Range("D92").NumberFormat = "@"
$ws.Range("D92").Value = "152445.71"
$ws.Range("D92").Style = "Normal"

$ws.Range("C93").NumberFormat = "@"
$ws.Range("C93").Value = "108"
$ws.Range("C93").Style = "Normal"
$ws.Range("D93").NumberFormat = "@"
$ws.Range("D93").Value = "265405.00"
$ws.Range("D93").Style = "Normal"

$ws.Range("C94").NumberFormat = "@"
$ws.Range("C94").Value = "22"
$ws.Range("C94").Style = "Normal"
$ws.Range("D94").NumberFormat = "@"
$ws.Range("D94").Value = "45500.00"
$ws.Range("D94").Style = "Normal"

$ws.Range("C95").NumberFormat = "@"
$ws.Range("C95").Value = "134"
$ws.Range("C95").Style = "Normal"
$ws.Range("D95").NumberFormat = "@"
$ws.Range("D95").Value = "387467.00"
$ws.Range("D95").Style = "Normal"

$ws.Range("C97").NumberFormat = "@"
$ws.Range("C97").Value = "8"
$ws.Range("C97").Style = "Normal"
$ws.Range("D97").NumberFormat = "@"
$ws.Range("D97").Value = "17500.00"
$ws.Range("D97").Style = "Normal"

$ws.Range("C98").NumberFormat = "@"
$ws.Range("C98").Value = "11"
$ws.Range("C98").Style = "Normal"
$ws.Range("D98").NumberFormat = "@"
$ws.Range("D98").Value = "25000.00"
$ws.Range("D98").Style = "Normal"

$ws.Range("C99").NumberFormat = "@"
$ws.Range("C99").Value = "44"
$ws.Range("C99").Style = "Normal"
$ws.Range("D99").NumberFormat = "@"
$ws.Range("D99").Value = "130500.00"
$ws.Range("D99").Style = "Normal"

$ws.Range("C100").NumberFormat = "@"
$ws.Range("C100").Value = "61"
$ws.Range("C100").Style = "Normal"
$ws.Range("D100").NumberFormat = "@"
$ws.Range("D100").Value = "147979.00"
$ws.Range("D100").Style = "Normal"

$ws.Range("C103").NumberFormat = "@"
$ws.Range("C103").Value = "20"
$ws.Range("C103").Style = "Normal"
$ws.Range("D103").NumberFormat = "@"
$ws.Range("D103").Value = "50830.00"
$ws.Range("D103").Style = "Normal"

$ws.Range("C104").NumberFormat = "@"
$ws.Range("C104").Value = "51"
$ws.Range("C104").Style = "Normal"
$ws.Range("D104").NumberFormat = "@"
$ws.Range("D104").Value = "103500.00"
$ws.Range("D104").Style = "Normal"

$ws.Range("C121").NumberFormat = "@"
$ws.Range("C121").Value = "61"
$ws.Range("C121").Style = "Normal"
$ws.Range("D121").NumberFormat = "@"
$ws.Range("D121").Value = "168877.00"
$ws.Range("D121").Style = "Normal"

$ws.Range("C122").NumberFormat = "@"
$ws.Range("C122").Value = "222"
$ws.Range("C122").Style = "Normal"
$ws.Range("D122").NumberFormat = "@"
$ws.Range("D122").Value = "603708.00"
$ws.Range("D122").Style = "Normal"

$ws.Range("C123").NumberFormat = "@"
$ws.Range("C123").Value = "76"
$ws.Range("C123").Style = "Normal"
$ws.Range("D123").NumberFormat = "@"
$ws.Range("D123").Value = "215531.45"
$ws.Range("D123").Style = "Normal"

$ws.Range("C124").NumberFormat = "@"
$ws.Range("C124").Value = "423"
$ws.Range("C124").Style = "Normal"
$ws.Range("D124").NumberFormat = "@"
$ws.Range("D124").Value = "1748149.06"
$ws.Range("D124").Style = "Normal"

$ws.Range("C128").NumberFormat = "@"
$ws.Range("C128").Value = "82"
$ws.Range("C128").Style = "Normal"
$ws.Range("D128").NumberFormat = "@"
$ws.Range("D128").Value = "244743.68"
$ws.Range("D128").Style = "Normal"

$ws.Range("C132").NumberFormat = "@"
$ws.Range("C132").Value = "72"
$ws.Range("C132").Style = "Normal"
$ws.Range("D132").NumberFormat = "@"
$ws.Range("D132").Value = "302086.75"
$ws.Range("D132").Style = "Normal"

$ws.Range("C134").NumberFormat = "@"
$ws.Range("C134").Value = "15"
$ws.Range("C134").Style = "Normal"
$ws.Range("D134").NumberFormat = "@"
$ws.Range("D134").Value = "44780.00"
$ws.Range("D134").Style = "Normal"

$ws.Range("C135").NumberFormat = "@"
$ws.Range("C135").Value = "203"
$ws.Range("C135").Style = "Normal"
$ws.Range("D135").NumberFormat = "@"
$ws.Range("D135").Value = "552620.00"
$ws.Range("D135").Style = "Normal"

$ws.Range("C136").NumberFormat = "@"
$ws.Range("C136").Value = "8"
$ws.Range("C136").Style = "Normal"
$ws.Range("D136").NumberFormat = "@"
$ws.Range("D136").Value = "16000.00"
$ws.Range("D136").Style = "Normal"

$ws.Range("C137").NumberFormat = "@"
$ws.Range("C137").Value = "17"
$ws.Range("C137").Style = "Normal"
$ws.Range("D137").NumberFormat = "@"
$ws.Range("D137").Value = "38500.00"
$ws.Range("D137").Style = "Normal"

$ws.Range("C138").NumberFormat = "@"
$ws.Range("C138").Value = "550"
$ws.Range("C138").Style = "Normal"
$ws.Range("D138").NumberFormat = "@"
$ws.Range("D138").Value = "1363646.00"
$ws.Range("D138").Style = "Normal"

$ws.Range("C139").NumberFormat = "@"
$ws.Range("C139").Value = "1748"
$ws.Range("C139").Style = "Normal"
$ws.Range("D139").NumberFormat = "@"
$ws.Range("D139").Value = "4677039.93"
$ws.Range("D139").Style = "Normal"

$ws.Range("C140").NumberFormat = "@"
$ws.Range("C140").Value = "2257"
$ws.Range("C140").Style = "Normal"
$ws.Range("D140").NumberFormat = "@"
$ws.Range("D140").Value = "5633837.29"
$ws.Range("D140").Style = "Normal"

$ws.Range("C141").NumberFormat = "@"
$ws.Range("C141").Value = "2419"
$ws.Range("C141").Style = "Normal"
$ws.Range("D141").NumberFormat = "@"
$ws.Range("D141").Value = "10041623.38"
$ws.Range("D141").Style = "Normal"

$ws.Range("C142").NumberFormat = "@"
$ws.Range("C142").Value = "341"
$ws.Range("C142").Style = "Normal"
$ws.Range("D142").NumberFormat = "@"
$ws.Range("D142").Value = "962007.04"
$ws.Range("D142").Style = "Normal"

$ws.Range("C143").NumberFormat = "@"
$ws.Range("C143").Value = "119"
$ws.Range("C143").Style = "Normal"
$ws.Range("D143").NumberFormat = "@"
$ws.Range("D143").Value = "294500.00"
$ws.Range("D143").Style = "Normal"

$ws.Range("C144").NumberFormat = "@"
$ws.Range("C144").Value = "234"
$ws.Range("C144").Style = "Normal"
$ws.Range("D144").NumberFormat = "@"
$ws.Range("D144").Value = "598800.00"
$ws.Range("D144").Style = "Normal"

$ws.Range("C145").NumberFormat = "@"
$ws.Range("C145").Value = "994"
$ws.Range("C145").Style = "Normal"
$ws.Range("D145").NumberFormat = "@"
$ws.Range("D145").Value = "2579433.25"
$ws.Range("D145").Style = "Normal"

$ws.Range("C146").NumberFormat = "@"
$ws.Range("C146").Value = "465"
$ws.Range("C146").Style = "Normal"
$ws.Range("D146").NumberFormat = "@"
$ws.Range("D146").Value = "1343454.49"
$ws.Range("D146").Style = "Normal"

$ws.Range("C147").NumberFormat = "@"
$ws.Range("C147").Value = "357"
$ws.Range("C147").Style = "Normal"
$ws.Range("D147").NumberFormat = "@"
$ws.Range("D147").Value = "893700.16"
$ws.Range("D147").Style = "Normal"

$ws.Range("C148").NumberFormat = "@"
$ws.Range("C148").Value = "143"
$ws.Range("C148").Style = "Normal"
$ws.Range("D148").NumberFormat = "@"
$ws.Range("D148").Value = "353000.00"
$ws.Range("D148").Style = "Normal"

$ws.Range("C149").NumberFormat = "@"
$ws.Range("C149").Value = "383"
$ws.Range("C149").Style = "Normal"
$ws.Range("D149").NumberFormat = "@"
$ws.Range("D149").Value = "1192090.28"
$ws.Range("D149").Style = "Normal"

$ws.Range("C150").NumberFormat = "@"
$ws.Range("C150").Value = "826"
$ws.Range("C150").Style = "Normal"
$ws.Range("D150").NumberFormat = "@"
$ws.Range("D150").Value = "1980242.82"
$ws.Range("D150").Style = "Normal"

$ws.Range("C192").NumberFormat = "@"
$ws.Range("C192").Value = "53"
$ws.Range("C192").Style = "Normal"
$ws.Range("D192").NumberFormat = "@"
$ws.Range("D192").Value = "133653.50"
$ws.Range("D192").Style = "Normal"

$ws.Range("C193").NumberFormat = "@"
$ws.Range("C193").Value = "110"
$ws.Range("C193").Style = "Normal"
$ws.Range("D193").NumberFormat = "@"
$ws.Range("D193").Value = "292000.00"
$ws.Range("D193").Style = "Normal"

$ws.Range("C194").NumberFormat = "@"
$ws.Range("C194").Value = "344"
$ws.Range("C194").Style = "Normal"
$ws.Range("D194").NumberFormat = "@"
$ws.Range("D194").Value = "930788.00"
$ws.Range("D194").Style = "Normal"

$ws.Range("C196").NumberFormat = "@"
$ws.Range("C196").Value = "614"
$ws.Range("C196").Style = "Normal"
$ws.Range("D196").NumberFormat = "@"
$ws.Range("D196").Value = "2241053.66"
$ws.Range("D196").Style = "Normal"

$ws.Range("C200").NumberFormat = "@"
$ws.Range("C200").Value = "151"
$ws.Range("C200").Style = "Normal"
$ws.Range("D200").NumberFormat = "@"
$ws.Range("D200").Value = "444633.00"
$ws.Range("D200").Style = "Normal"

$ws.Range("C202").NumberFormat = "@"
$ws.Range("C202").Value = "74"
$ws.Range("C202").Style = "Normal"
$ws.Range("D202").NumberFormat = "@"
$ws.Range("D202").Value = "175005.00"
$ws.Range("D202").Style = "Normal"

$ws.Range("C204").NumberFormat = "@"
$ws.Range("C204").Value = "112"
$ws.Range("C204").Style = "Normal"
$ws.Range("D204").NumberFormat = "@"
$ws.Range("D204").Value = "502180.50"
$ws.Range("D204").Style = "Normal"

$ws.Range("C205").NumberFormat = "@"
$ws.Range("C205").Value = "130"
$ws.Range("C205").Style = "Normal"
$ws.Range("D205").NumberFormat = "@"
$ws.Range("D205").Value = "288696.77"
$ws.Range("D205").Style = "Normal"

$ws.Range("C223").NumberFormat = "@"
$ws.Range("C223").Value = "58"
$ws.Range("C223").Style = "Normal"
$ws.Range("D223").NumberFormat = "@"
$ws.Range("D223").Value = "153500.00"
$ws.Range("D223").Style = "Normal"

$ws.Range("C224").NumberFormat = "@"
$ws.Range("C224").Value = "158"
$ws.Range("C224").Style = "Normal"
$ws.Range("D224").NumberFormat = "@"
$ws.Range("D224").Value = "431905.00"
$ws.Range("D224").Style = "Normal"

$ws.Range("C231").NumberFormat = "@"
$ws.Range("C231").Value = "28"
$ws.Range("C231").Style = "Normal"
$ws.Range("D231").NumberFormat = "@"
$ws.Range("D231").Value = "84587.00"
$ws.Range("D231").Style = "Normal"

